# Apply updated cryptocurrency price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.021.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "'1.901.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'0.7372"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'241.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.3061"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("D9").Value = "'26.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("D10").Value = "'0.06883"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").Value = "'0.08046"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'0.7617"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "'1.910.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'5.222"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "'91.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "'30.030.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "'6.055"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'14.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "'0.000007735"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("D20").Value = "'237.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.152.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'7.051"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("D25").Value = "'9.279"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "'166.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'18.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "'0.1260"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").Value = "'2.027"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.37%  "
$ws.Range("D30").Value = "'1.363"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'4.269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "'4.034"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "'0.05462"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("D35").Value = "'1.290"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").Value = "'0.7343"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "'0.01940"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'6.270"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").Value = "'0.4429"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").Value = "'73.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("D43").Value = "'1.956"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'0.8340"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'101.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.600"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'9.835"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").Value = "'2.058.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'36.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05966"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
